$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.899.05"
$ws.Range("E2").Value = "  -1.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.476.15"
$ws.Range("E3").Value = "  -2.45%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.11"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.73"
$ws.Range("E6").Value = "  -4.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.474.72"
$ws.Range("E7").Value = "  -2.47%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -2.80%  "

$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.74"
$ws.Range("E11").Value = "  +3.75%  "

$ws.Range("E12").Value = "  -3.66%  "

$ws.Range("E13").Value = "  -4.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.062.77"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.08"
$ws.Range("E15").Value = "  -6.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.471.38"
$ws.Range("E16").Value = "  -2.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.862.32"
$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  -5.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.14"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.17"
$ws.Range("E21").Value = "  -5.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.12"
$ws.Range("E22").Value = "  -4.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.605"
$ws.Range("E23").Value = "  -6.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.15"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.611.94"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("E27").Value = "  -9.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  -6.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.32"
$ws.Range("E29").Value = "  -10.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("E31").Value = "  -6.91%  "

$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.32"
$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.465.29"
$ws.Range("E35").Value = "  -2.57%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -6.30%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.90"
$ws.Range("E37").Value = "  -7.68%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.89"
$ws.Range("E39").Value = "  -4.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.79"
$ws.Range("E41").Value = "  -4.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0884"
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.40"
$ws.Range("E43").Value = "  -4.34%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  -13.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.897"
$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.00"
$ws.Range("E46").Value = "  -6.97%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.28"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("E48").Value = "  -7.71%  "

$ws.Range("E49").Value = "  -4.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.39"
$ws.Range("E50").Value = "  -9.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.968"
$ws.Range("E51").Value = "  -4.93%  "
